$d = $word.ActiveDocument

# Locate the "Abstract" paragraph that begins "Federizer is a digital data
# transfer system akin to email system, ..." -- do this by scanning
# paragraphs rather than hard-coding an index, so we are robust to any
# paragraph numbering differences.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Federizer is a digital data transfer system")) {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)

# --- "data" -> "media" -------------------------------------------------
$p.Range.Find.Execute("data", $false, $false, $false, $false, $false, $true, 1, $false, "media", 2)

# --- "transfer" -> "exchange and storage" ------------------------------
$p.Range.Find.Execute("transfer", $false, $false, $false, $false, $false, $true, 1, $false, "exchange and storage", 2)

# --- Rewrite the trailing sentence(s) -----------------------------------
# Old:  " akin to email system, aligned with emerging and future business
#        needs. It is a replacement for the current email system that no
#        longer meets security standards and functional requirements."
# New:  " akin to email system. It is aligned with emerging and future
#        business needs. It is a replacement for the current email system
#        that no longer meets security standards and functional
#        requirements."
# ... with the two "It is" fragments recolored to C9211E.
$oldPhrase = " akin to email system, aligned with emerging and future business needs. It is a replacement for the current email system that no longer meets security standards and functional requirements."
$full = $p.Range.Text
$phraseStart = $p.Range.Start + $full.IndexOf($oldPhrase)
$phraseEnd = $phraseStart + $oldPhrase.Length
$phraseRange = $d.Range($phraseStart, $phraseEnd)

$newPhrase = " akin to email system. It is aligned with emerging and future business needs. It is a replacement for the current email system that no longer meets security standards and functional requirements."
$phraseRange.Text = $newPhrase

$pStart = $p.Range.Start
$full2 = $p.Range.Text

$marker1 = "akin to email system. "
$idx1 = $full2.IndexOf($marker1) + $marker1.Length
$r1 = $d.Range($pStart + $idx1, $pStart + $idx1 + 5)
$r1.Font.Color = 1974729

$marker2 = "needs. "
$idx2 = $full2.IndexOf($marker2) + $marker2.Length
$r2 = $d.Range($pStart + $idx2, $pStart + $idx2 + 6)
$r2.Font.Color = 1974729
